$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Rename header row (row 1) columns: "_old" suffix -> "_FV2404", "_new" suffix -> "_FV2410" ---
$newHeadersFV2404 = @(
    "Segmentname_FV2404", "Segmentgruppe_FV2404", "Segment_FV2404", "Datenelement_FV2404",
    "Segment ID_FV2404", "Code_FV2404", "Qualifier_FV2404", "Beschreibung_FV2404",
    "Bedingungsausdruck_FV2404", "Bedingung_FV2404"
)
$newHeadersFV2410 = @(
    "Segmentname_FV2410", "Segmentgruppe_FV2410", "Segment_FV2410", "Datenelement_FV2410",
    "Segment ID_FV2410", "Code_FV2410", "Qualifier_FV2410", "Beschreibung_FV2410",
    "Bedingungsausdruck_FV2410", "Bedingung_FV2410"
)

for ($i = 0; $i -lt 10; $i++) {
    # Columns A-J (1-10)
    $ws.Cells.Item(1, $i + 1).Value = $newHeadersFV2404[$i]
    # Columns L-U (12-21)
    $ws.Cells.Item(1, $i + 12).Value = $newHeadersFV2410[$i]
}

# --- Add an Excel Table (ListObject) over the full used range, with autofilter ---
$rangeAddress = "A1:U56"
$listObj = $ws.ListObjects.Add(1, $ws.Range($rangeAddress), $null, 1)
$listObj.Name = "Table1"

# --- Freeze the header row (row 1) ---
$ws.Activate()
[void]$ws.Range("A2").Select()
$excel.ActiveWindow.FreezePanes = $true
